$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = "tWlJzrJaXYzgyOCBOLtu"
$ws.Range("F2").Value = "'-179.958101"
$ws.Range("G2").Value = "'-2.5592605"

# Row 3
$ws.Range("E3").Value = "gPcUrDKYsImEbbgIsjvI"
$ws.Range("F3").Value = "'106.989965"
$ws.Range("G3").Value = "'64.934071"

# Row 4 (only the text column changes; longitude/latitude stay blank)
$ws.Range("E4").Value = "PoSeWuySWgMXFWTtqxet"

# Row 5
$ws.Range("E5").Value = "wKGNGWPcXSdbsqfXZicm"
$ws.Range("F5").Value = "'-122.893621"
$ws.Range("G5").Value = "'1.081778"

# Row 6
$ws.Range("E6").Value = "GpbjzsVxKUVTCvGBLPUF"
$ws.Range("F6").Value = "'-83.086954"
$ws.Range("G6").Value = "'78.3302665"
